$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Text number format first so numeric-looking strings (e.g. "319.15")
# are preserved exactly as text, matching the source data (avoids float drift).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.268.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.424.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.15'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.55'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.48'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.26'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.804.40'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.405.88'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.163.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.96'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.64'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.70'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.64'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.91'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.11'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.19%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.60%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0762'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.52'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.61'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.935.55'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.95'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.05%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.56'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.82'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.86%  '
